{"js": "// Replace the 25 two-digit multiplication problems in the table with\n// their new values (order-preserving 1:1 text swap, per the diff).\nconst replacements = [\n  [\"23\u00d767=\", \"72\u00d770=\"],\n  [\"41\u00d786=\", \"85\u00d744=\"],\n  [\"11\u00d797=\", \"49\u00d721=\"],\n  [\"40\u00d774=\", \"49\u00d730=\"],\n  [\"97\u00d795=\", \"32\u00d769=\"],\n  [\"31\u00d717=\", \"41\u00d734=\"],\n  [\"56\u00d784=\", \"53\u00d713=\"],\n  [\"22\u00d787=\", \"49\u00d737=\"],\n  [\"47\u00d748=\", \"12\u00d743=\"],\n  [\"22\u00d768=\", \"94\u00d742=\"],\n  [\"12\u00d749=\", \"64\u00d754=\"],\n  [\"23\u00d725=\", \"84\u00d730=\"],\n  [\"77\u00d722=\", \"35\u00d754=\"],\n  [\"51\u00d715=\", \"42\u00d741=\"],\n  [\"18\u00d753=\", \"97\u00d749=\"],\n  [\"94\u00d733=\", \"47\u00d755=\"],\n  [\"14\u00d761=\", \"40\u00d772=\"],\n  [\"87\u00d728=\", \"76\u00d751=\"],\n  [\"39\u00d789=\", \"39\u00d785=\"],\n  [\"76\u00d791=\", \"98\u00d780=\"],\n  [\"40\u00d798=\", \"49\u00d731=\"],\n  [\"65\u00d724=\", \"17\u00d711=\"],\n  [\"41\u00d797=\", \"43\u00d770=\"],\n  [\"28\u00d745=\", \"59\u00d733=\"],\n  [\"21\u00d732=\", \"41\u00d740=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication problems in the table with\n# their new values (order-preserving 1:1 text swap, per the diff).\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"23\u00d767=\", \"72\u00d770=\"),\n    @(\"41\u00d786=\", \"85\u00d744=\"),\n    @(\"11\u00d797=\", \"49\u00d721=\"),\n    @(\"40\u00d774=\", \"49\u00d730=\"),\n    @(\"97\u00d795=\", \"32\u00d769=\"),\n    @(\"31\u00d717=\", \"41\u00d734=\"),\n    @(\"56\u00d784=\", \"53\u00d713=\"),\n    @(\"22\u00d787=\", \"49\u00d737=\"),\n    @(\"47\u00d748=\", \"12\u00d743=\"),\n    @(\"22\u00d768=\", \"94\u00d742=\"),\n    @(\"12\u00d749=\", \"64\u00d754=\"),\n    @(\"23\u00d725=\", \"84\u00d730=\"),\n    @(\"77\u00d722=\", \"35\u00d754=\"),\n    @(\"51\u00d715=\", \"42\u00d741=\"),\n    @(\"18\u00d753=\", \"97\u00d749=\"),\n    @(\"94\u00d733=\", \"47\u00d755=\"),\n    @(\"14\u00d761=\", \"40\u00d772=\"),\n    @(\"87\u00d728=\", \"76\u00d751=\"),\n    @(\"39\u00d789=\", \"39\u00d785=\"),\n    @(\"76\u00d791=\", \"98\u00d780=\"),\n    @(\"40\u00d798=\", \"49\u00d731=\"),\n    @(\"65\u00d724=\", \"17\u00d711=\"),\n    @(\"41\u00d797=\", \"43\u00d770=\"),\n    @(\"28\u00d745=\", \"59\u00d733=\"),\n    @(\"21\u00d732=\", \"41\u00d740=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
